$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Unique" column (L) to the EXPERIMENT_TYPE property-assignments
# table: header in L4 (mirrors K4's style) and a "FALSE" value for every
# property row (L5:L8, mirroring K5:K8's style).

$ws.Range("L4").Value = "Unique"
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)

# Set L5:L8 to the literal text "FALSE" (not the boolean) by seeding the
# value from an already-text "FALSE" cell (B5), which keeps it a string.
# (Done per-cell: pasting one source cell onto a multi-cell target only
# fills the first cell, unlike a real Excel "spill" paste.)
$ws.Range("B5").Copy()
$ws.Range("L5").PasteSpecial(-4163)
$ws.Range("B5").Copy()
$ws.Range("L6").PasteSpecial(-4163)
$ws.Range("B5").Copy()
$ws.Range("L7").PasteSpecial(-4163)
$ws.Range("B5").Copy()
$ws.Range("L8").PasteSpecial(-4163)

# Copy the K column's number format / style onto L, row by row.
$ws.Range("K5").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("K6").Copy()
$ws.Range("L6").PasteSpecial(-4122)
$ws.Range("K7").Copy()
$ws.Range("L7").PasteSpecial(-4122)
$ws.Range("K8").Copy()
$ws.Range("L8").PasteSpecial(-4122)

$ws.Range("L7:L8").Select()
